# Item.xlsx edit: add a new "Equip_Weapon_1" row (开山斧 / 开山斧武器) to the
# XML-mapped table on Sheet1, growing the table/used range from A1:K8 to
# A1:K9, and move the active selection to K13.
#
# Note: the source absPath (xl/workbook.xml mc:Choice/x15ac:absPath) that
# recorded the author's local save folder is not exposed anywhere on the
# Excel COM object model (Workbook.Path/.FullName are read-only, and it is
# preserved verbatim on save regardless of what is touched) so it cannot be
# changed from a COM/VBA-level script.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the new data row (row 9) -------------------------------------
# Columns: ID, ItemType, ItemSubType, Level, ShowName, Desc, Icon,
#          CoolDownTime, OverlayCount, BuyPrice, SalePrice
# Write the text-valued columns first, in left-to-right order, so the new
# shared-string entries land in the same order as the reference edit
# (Equip_Weapon_1, 开山斧, 开山斧武器, 50004).
$ws.Cells.Item(9, 1).Value = "Equip_Weapon_1"
$ws.Cells.Item(9, 5).Value = "开山斧"
$ws.Cells.Item(9, 6).Value = "开山斧武器"

# Icon is stored as text (matches style used by the other rows' Icon cells);
# set the text number format before assigning so "50004" is kept as a
# string instead of being coerced to a number.
$ws.Cells.Item(9, 7).NumberFormat = "@"
$ws.Cells.Item(9, 7).Value = "50004"

# Numeric columns.
$ws.Cells.Item(9, 2).Value = 1
$ws.Cells.Item(9, 3).Value = 1
$ws.Cells.Item(9, 4).Value = 1
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 10000
$ws.Cells.Item(9, 10).Value = 100
$ws.Cells.Item(9, 11).Value = 100

# --- Grow the XML-mapped table to include the new row ---------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:K9")) | Out-Null

# --- Move the selection cursor to K13, as in the edited workbook ----------
$ws.Range("K13").Select() | Out-Null
